$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column formatting -----------------------------------------------------
# Columns F:G get a dedicated "x"/"y" number column style (center aligned)
# distinct from the default column style used everywhere else (H onward).
$ws.Columns("F:G").ColumnWidth = 10.6

# --- Header cells (F6 = "x", G6 = "y") get center alignment ---------------
$ws.Range("F6").HorizontalAlignment = -4108
$ws.Range("G6").HorizontalAlignment = -4108

# --- Row 8 : vdda supply -----------------------------------------------
$ws.Range("F8").Formula = "=10"
$ws.Range("G8").Value = 10
$ws.Range("F8").HorizontalAlignment = -4108
$ws.Range("G8").HorizontalAlignment = -4108

# --- Row 9 : vdd supply -------------------------------------------------
$ws.Range("F9").Formula = "=F8+10"
$ws.Range("G9").Formula = "=`$G`$8"
$ws.Range("F9").HorizontalAlignment = -4108
$ws.Range("G9").HorizontalAlignment = -4108

# --- Row 10 : gnd ground -------------------------------------------------
$ws.Range("F10").Formula = "=F9+10"
$ws.Range("G10").Formula = "=`$G`$8"
$ws.Range("F10").HorizontalAlignment = -4108
$ws.Range("G10").HorizontalAlignment = -4108

# --- Row 11 (new): intermediate pin x position, no signal row here --------
$ws.Range("F11").Formula = "=F10+10"
$ws.Range("F11").HorizontalAlignment = -4108

# --- Row 12 : raddr sig_in -------------------------------------------------
$ws.Range("F12").Formula = "=F11+10"
$ws.Range("G12").Formula = "=`$G`$8"
$ws.Range("F12").HorizontalAlignment = -4108
$ws.Range("G12").HorizontalAlignment = -4108

# --- Row 13 : rdata sig_out -------------------------------------------------
$ws.Range("F13").Formula = "=F12+10"
$ws.Range("G13").Formula = "=`$G`$8"
$ws.Range("F13").HorizontalAlignment = -4108
$ws.Range("G13").HorizontalAlignment = -4108

# --- Row 14 : rclk clk_in -------------------------------------------------
$ws.Range("F14").Formula = "=F13+10"
$ws.Range("G14").Formula = "=`$G`$8"
$ws.Range("F14").HorizontalAlignment = -4108
$ws.Range("G14").HorizontalAlignment = -4108

# --- Row 15 (new): intermediate pin x position, no signal row here --------
$ws.Range("F15").Formula = "=F14+10"
$ws.Range("F15").HorizontalAlignment = -4108

# --- Row 16 : waddr sig_in -------------------------------------------------
$ws.Range("F16").Formula = "=F15+10"
$ws.Range("G16").Formula = "=`$G`$8"
$ws.Range("F16").HorizontalAlignment = -4108
$ws.Range("G16").HorizontalAlignment = -4108

# --- Row 17 : wdata sig_in -------------------------------------------------
$ws.Range("F17").Formula = "=F16+10"
$ws.Range("G17").Formula = "=`$G`$8"
$ws.Range("F17").HorizontalAlignment = -4108
$ws.Range("G17").HorizontalAlignment = -4108

# --- Row 18 : wen sig_in -------------------------------------------------
$ws.Range("F18").Formula = "=F17+10"
$ws.Range("G18").Formula = "=`$G`$8"
$ws.Range("F18").HorizontalAlignment = -4108
$ws.Range("G18").HorizontalAlignment = -4108

# --- Row 19 : wclk clk_in -------------------------------------------------
$ws.Range("F19").Formula = "=F18+10"
$ws.Range("G19").Formula = "=`$G`$8"
$ws.Range("F19").HorizontalAlignment = -4108
$ws.Range("G19").HorizontalAlignment = -4108

# --- Move the active selection to reflect where editing left off ----------
[void]$ws.Range("G22").Select()
